$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "usproj_category_longname"
$ws.Range("B2").Value = "Abandoned Coal Mines"
$ws.Range("B3").Value = "Abandoned Oil Wells"

$newRows = @(
    @("CroplandConvertCarbonStock", "Land Converted to Cropland: Carbon Stocks"),
    @("CroplandRemainCarbonStock", "Cropland Remaining Cropland: Carbon Stocks"),
    @("ForestLandConvertCarbonStock", "Land Converted to Forest Land: Carbon Stock"),
    @("ForestLandDrainedOrg", "Forest Land Remaining Forest Land: Drained Organic Soils"),
    @("ForestLandFires", "Forest Land Remaining Forest Land: Forest Fires"),
    @("ForestLandRemainCarbonStock", "Forest Land Remaining Forest Land: Carbon Stocks"),
    @("ForestLandSoils", "Forest Land Remaining Forest Land: Forest Soils"),
    @("GrasslandConvertCarbonStock", "Land Converted to Grassland: Carbon Stocks"),
    @("GrasslandFires", "Grassland Remaining Grassland: Grassland Fires"),
    @("GrasslandRemainCarbonStock", "Grassland Remaining Grassland: Carbon Stocks"),
    @("SettlementsConvertCarbonStock", "Land Converted to Settlements: Carbon Stocks"),
    @("SettlementsRemainCarbonStock", "Settlements Remaining Settlements: Carbon Stocks"),
    @("SettlementsSoils", "Settlements Remaining Settlements: Settlement Soils"),
    @("WetlandsCoastal", "Coastal Wetlands Remaining Coastal Wetlands"),
    @("WetlandsConvertCarbonStock", "Land Converted to Wetlands: Carbon Stocks"),
    @("WetlandsConvertCoastal", "Lands Converted to Coastal Wetlands"),
    @("WetlandsConvertFlooded", "Lands Converted to Flooded Lands"),
    @("WetlandsFlooded", "Flooded Land Remaining Flooded Land"),
    @("WetlandsPeatlands", "Peatlands Remaining Peatlands"),
    @("WetlandsRemainCarbonStock", "Wetlands Remaining Wetlands: Carbon Stocks")
)

$r = 57
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

Write-Output "done"
